$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.219.90"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "1.604.60"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3771"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.98%  "

$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("E10").Value = "  +1.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08145"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.602"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.357"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001246"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Value = "1.600.53"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06937"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.534"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").Value = "23.225.64"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.432"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.062"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.50"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.293"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.418"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.765"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").Value = "1.778.78"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9559"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02765"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07440"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2516"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08775"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.409"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7100"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6531"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.330"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.011"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07950"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.199"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.194"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.12%  "
